$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "30-10-2025"
$ws.Range("B46").Value = "The price of gold in India today is ₹12,049 per gram for 24 karat gold, ₹11,045 per gram for 22 karat gold and ₹9,037 per gram for 18 karat gold (also called 999 gold)."
